$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (20 questions)
$ws.Range("C4").Value = "0.41/s"
$ws.Range("D4").Value = "2.42/s"
$ws.Range("E4").Value = "0.59/s"
$ws.Range("F4").Value = "1.76/s"
$ws.Range("G4").Value = "5.18/s"

# Row 5 (40 questions)
$ws.Range("C5").Value = "0.37/s"
$ws.Range("D5").Value = "2.48/s"
$ws.Range("E5").Value = "0.59/s"
$ws.Range("F5").Value = "1.58/s"
$ws.Range("G5").Value = "5.02/s"
$ws.Range("H5").Value = 0.84

# Row 6 (60 questions)
$ws.Range("C6").Value = "0.53/s"
$ws.Range("D6").Value = "2.68/s"
$ws.Range("E6").Value = "0.59/s"
$ws.Range("F6").Value = "1.83/s"
$ws.Range("G6").Value = "5.63/s"
$ws.Range("H6").Value = 0.85

# Update the active selection to match the authored change
$ws.Range("H9").Select()
